# Update loading_percent values for the 380 kV case (rows 2-25, cols C:O excluding B/H/J/N which stay 0)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colIndexes = @(3,4,5,6,7,9,11,12,13,15)

$data = @(
    @(9.810751046378549,5.727113849605437,12.4771394701916,27.16965059071424,3.655525909044426,25.19632950114741,14.91304794498668,9.555184319125832,16.96766325713664,24.42635460740105),
    @(9.779030120599835,5.686869527530178,12.5049828859166,27.28401472622143,3.657680072809946,25.31508253938248,14.36099799845857,9.586809169504368,16.73738966162511,24.56239591593863),
    @(9.761197709895399,5.661923362733416,12.52433206592875,27.3630440493674,3.659072095114931,25.39464822886116,14.01122371367977,9.607466228754946,16.59602619840617,24.65256077655153),
    @(9.754349513768176,5.651703027520118,12.53278312648029,27.3974534928631,3.659656851598189,25.42873919026896,13.86615798460755,9.61619618658251,16.53848404602807,24.69096734942185),
    @(9.753237804716813,5.650002791989272,12.53422059525673,27.40329996885592,3.65975500832609,25.43450051890177,13.84192307002534,9.617664648045857,16.52893482252854,24.69744507464632),
    @(9.761103651078175,5.661785742386609,12.52444374839569,27.3634991964168,3.659079910430975,25.39510124747915,14.00927728462222,9.607582700037257,16.59524982618393,24.65307201120392),
    @(9.799475603335228,5.713289445960248,12.48627195920623,27.20724859082581,3.656254303255698,25.23589210755749,14.72505264806462,9.565831604720676,16.88829446205846,24.47188170172155),
    @(9.887525880579483,5.812207779421975,12.42931472337027,26.97120546515423,3.651261066259672,24.97668150826986,16.03559359851403,9.493772174852436,17.46062242867156,24.16942693375184),
    @(9.95966993065108,5.883339612954066,12.3984024970275,26.84129068187256,3.647922897579864,24.8188554301491,16.93326400341038,9.446787154650846,17.87635386011518,23.97973167773438),
    @(9.994026059008515,5.915306272928659,12.38671753855376,26.79175689315633,3.646475252347972,24.75420329984151,17.32604882252921,9.426700573081456,18.06376352505727,23.90055971176758),
    @(10.00724938959442,5.927350040532585,12.38263472511215,26.7743842898924,3.645937205106452,24.73075379101618,17.47244809343837,9.419279017227149,18.13442985839614,23.8716089025385),
    @(10.00439214269771,5.924759005632674,12.38349881818659,26.77806406094943,3.646052632837233,24.73575803794387,17.44102391673268,9.420869167743664,18.11922488471155,23.8777980786947),
    @(9.995109718349115,5.916298378190751,12.38637478752158,26.79029983458615,3.646430783858231,24.75225336307971,17.3381406014064,9.426086295256136,18.06958364568386,23.89815723692643),
    @(9.989451524781007,5.911107867333961,12.38818094796264,26.79797520318607,3.646663731882826,24.76249189797463,17.27481412671951,9.429305993959932,18.0391360436822,23.91076210658255),
    @(9.957454965458464,5.881242207931966,12.39921398975576,26.84472111574534,3.648018926978494,24.82322482885384,16.90727234506141,9.448125735318239,17.86406699044263,23.98504955111252),
    @(9.938214530086537,5.862816723251443,12.40659141861149,26.8758548866916,3.648868418924314,24.86231620058802,16.67773331060423,9.46000049331411,17.75618977761428,24.03245059229298),
    @(9.927293131400591,5.852182311570682,12.41105846490959,26.89466199949165,3.649363701132835,24.88547259726144,16.54424811056843,9.466951708597861,17.69398406483456,24.0603843600416),
    @(9.923620498688873,5.848575556742283,12.41260935101753,26.90118401120534,3.649532543515516,24.89342821302842,16.49880478901497,9.469326090545806,17.67289694778388,24.06995713958987),
    @(9.940247735152177,5.864781967838508,12.40578292033724,26.87244744919663,3.648777298394573,24.85808525947248,16.70232005492674,9.458723866707702,17.76769021261053,24.02733529171223),
    @(9.997830462160406,5.918785173156847,12.38552076201104,26.78666823778684,3.646319436822808,24.74738021086651,17.36842416332169,9.424548885396076,18.08417310354713,23.89214926067648),
    @(10.03670425162079,5.953719735659502,12.37427188985297,26.73868194122507,3.644772189052233,24.6810508325459,17.79008601342642,9.403290566902992,18.28923023400749,23.80980314457209),
    @(10.01584571768385,5.935109047867324,12.38009316772294,26.76355131144196,3.645592592907103,24.71589916204026,17.56631834020932,9.414538079474045,18.17996847784419,23.85320134244851),
    @(9.93932808592057,5.86389360917954,12.40614773949994,26.87398512433541,3.64881847249538,24.85999594318231,16.69120911496814,9.459300642511332,17.76249143997626,24.0296457942443),
    @(9.86237007254149,5.785698954208434,12.44280438390547,27.02746358244865,3.652553597309529,25.04109914112944,15.69196395550025,9.512218194209657,17.30639080222418,24.24555990906726)
)

$startRow = 2
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $colIndexes.Length; $c++) {
        $ws.Cells.Item($startRow + $r, $colIndexes[$c]).Value = $rowVals[$c]
    }
}
